$d = $word.ActiveDocument

# --- Step 1: remove the "Version 0.2" heading and its (not yet implemented)
#     bullet items -- delete from the bottom up so earlier paragraph
#     indices stay stable.
$vHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Version 0.2") {
        $vHeading = $i
        break
    }
}

# Find the next Heading1 paragraph after it (start of "Version 0.1"); the
# block to delete is everything from the "Version 0.2" heading up to (but
# not including) that next heading.
$nextHeadingStart = $null
for ($i = $vHeading + 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Style.NameLocal -eq "Heading 1") {
        $nextHeadingStart = $i
        break
    }
}

$pStart = $d.Paragraphs.Item($vHeading)
$pEnd = $d.Paragraphs.Item($nextHeadingStart - 1)
$blockRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$blockRange.Delete()

# --- Step 2: remove the four "materials needed" bullet items that follow
#     the screws line (Bolts M6, Lock nuts, Bolts M3, Nuts M3).
$boltsStart = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "3 x Bolts M6 x 50") {
        $boltsStart = $i
        break
    }
}
$nutsEnd = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "4 x Nuts M3") {
        $nutsEnd = $i
        break
    }
}

$p13 = $d.Paragraphs.Item($boltsStart)
$p16 = $d.Paragraphs.Item($nutsEnd)
$rng2 = $d.Range($p13.Range.Start, $p16.Range.End)
$rng2.Delete()

# --- Step 3: clear the "Screws 3 x 10 mm" bullet's text and turn it into
#     an empty ListParagraph with a plain 360-twip left indent (no
#     numbering).
$screwIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Screws 3 x 10 mm") {
        $screwIdx = $i
        break
    }
}

$pScrew = $d.Paragraphs.Item($screwIdx)
$clearRng = $d.Range($pScrew.Range.Start, $pScrew.Range.End - 1)
$clearRng.Delete()

$pScrew2 = $d.Paragraphs.Item($screwIdx)
# Re-applying the style from "Normal" drops the inherited <w:numPr/>
# entirely, leaving a clean paragraph we can then give an explicit indent.
$pScrew2.Style = "Normal"
$pScrew2.Style = "List Paragraph"
$pScrew2.LeftIndent = 18
